# Updates the cryptocurrency price/volume table on Sheet1 to reflect the
# latest values pulled by the GitHub Actions refresh job.
# Cells in column D that contain plain numeric-looking text (e.g. "0.999")
# are written with a leading apostrophe so Excel keeps them as text
# (matching the original inline-string/text formatting) instead of
# silently converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.949.46"
$ws.Range("E2").Value = "  -3.25%  "
$ws.Range("D3").Value = "3.838.83"
$ws.Range("E3").Value = "  -2.61%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'600.24"
$ws.Range("E5").Value = "  -1.84%  "
$ws.Range("D6").Value = "'167.55"
$ws.Range("E6").Value = "  -2.49%  "
$ws.Range("D7").Value = "3.836.30"
$ws.Range("E7").Value = "  -2.67%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -2.05%  "
$ws.Range("E10").Value = "  -4.39%  "
$ws.Range("D11").Value = "'6.46"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "'0.458"
$ws.Range("E12").Value = "  -3.11%  "
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("D14").Value = "'36.95"
$ws.Range("E14").Value = "  -4.71%  "
$ws.Range("D15").Value = "4.478.65"
$ws.Range("E15").Value = "  -2.59%  "
$ws.Range("D16").Value = "3.842.23"
$ws.Range("E16").Value = "  -2.44%  "
$ws.Range("D17").Value = "68.072.22"
$ws.Range("E17").Value = "  -3.03%  "
$ws.Range("D18").Value = "'18.23"
$ws.Range("E18").Value = "  -2.06%  "
$ws.Range("D19").Value = "'7.39"
$ws.Range("E19").Value = "  -3.94%  "
$ws.Range("E20").Value = "  -0.95%  "
$ws.Range("D21").Value = "'11.03"
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("D22").Value = "'466.36"
$ws.Range("E22").Value = "  -6.41%  "
$ws.Range("E23").Value = "  -2.36%  "
$ws.Range("E24").Value = "  -4.01%  "
$ws.Range("D25").Value = "'82.89"
$ws.Range("E25").Value = "  -3.79%  "
$ws.Range("E26").Value = "  -3.55%  "
$ws.Range("E27").Value = "  -2.66%  "
$ws.Range("D28").Value = "'10.06"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").Value = "'2.97"
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("D31").Value = "3.984.45"
$ws.Range("E31").Value = "  -2.53%  "
$ws.Range("D32").Value = "'7.64"
$ws.Range("E32").Value = "  -3.46%  "
$ws.Range("E33").Value = "  -6.00%  "
$ws.Range("D34").Value = "'31.24"
$ws.Range("E34").Value = "  -3.73%  "
$ws.Range("D35").Value = "'9.55"
$ws.Range("E35").Value = "  -1.28%  "
$ws.Range("D36").Value = "3.796.55"
$ws.Range("E36").Value = "  -2.66%  "
$ws.Range("E37").Value = "  -4.04%  "
$ws.Range("D38").Value = "'3.62"
$ws.Range("E38").Value = "  +8.96%  "
$ws.Range("D39").Value = "'0.140"
$ws.Range("E39").Value = "  -1.23%  "
$ws.Range("E40").Value = "  -3.18%  "
$ws.Range("E41").Value = "  -4.55%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("E43").Value = "  -5.37%  "
$ws.Range("D44").Value = "'1.98"
$ws.Range("E44").Value = "  -7.68%  "
$ws.Range("D45").Value = "'421.57"
$ws.Range("E45").Value = "  -4.15%  "
$ws.Range("D46").Value = "'8.70"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("B48").Value = "FLOKI"
$ws.Range("C48").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D48").Value = "'0.000291"
$ws.Range("E48").Value = "  +4.15%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "'46.89"
$ws.Range("E49").Value = "  -2.99%  "
$ws.Range("D50").Value = "'142.17"
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0358"
$ws.Range("E51").Value = "  -3.25%  "
